# Fruta / hortaliza, semanal
# A new week of price observations is inserted at the top of the data block
# (rows 126-127), pushing all the existing rows down by two (the sheet grows
# from 152 to 154 used rows). The two new rows reuse the same static
# dimension values (Mercado, Region, Producto, etc.) as the block they are
# inserted into and only carry fresh Fecha / Volumen / Precio* / Origen data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 126; Excel shifts rows
# 126:152 down to 128:154 and the new rows inherit row 126's formatting
# (so column D keeps its date number format).
$ws.Rows("126:127").Insert()

# New row 126 - "Primera" quality
$ws.Range("A126").Value = 11
$ws.Range("B126").Value = "Vega Monumental Concepción"
$ws.Range("C126").Value = "Bíobío"
$ws.Range("D126").Value = 45258
$ws.Range("E126").Value = 8
$ws.Range("F126").Value = "Fruta"
$ws.Range("G126").Value = 100101
$ws.Range("H126").Value = "Berries"
$ws.Range("I126").Value = 100101001
$ws.Range("J126").Value = "Arándano (blue)"
$ws.Range("K126").Value = "Sin especificar"
$ws.Range("L126").Value = "Primera"
$ws.Range("M126").Value = 200
$ws.Range("N126").Value = 6000
$ws.Range("O126").Value = 6500
$ws.Range("P126").Value = 6250
$ws.Range("Q126").Value = "$/bandeja 2 kilos"
$ws.Range("R126").Value = "Provincia de Curicó"
$ws.Range("S126").Value = 3125
$ws.Range("T126").Value = 2

# New row 127 - "Segunda" quality
$ws.Range("A127").Value = 11
$ws.Range("B127").Value = "Vega Monumental Concepción"
$ws.Range("C127").Value = "Bíobío"
$ws.Range("D127").Value = 45258
$ws.Range("E127").Value = 8
$ws.Range("F127").Value = "Fruta"
$ws.Range("G127").Value = 100101
$ws.Range("H127").Value = "Berries"
$ws.Range("I127").Value = 100101001
$ws.Range("J127").Value = "Arándano (blue)"
$ws.Range("K127").Value = "Sin especificar"
$ws.Range("L127").Value = "Segunda"
$ws.Range("M127").Value = 100
$ws.Range("N127").Value = 5000
$ws.Range("O127").Value = 5000
$ws.Range("P127").Value = 5000
$ws.Range("Q127").Value = "$/bandeja 2 kilos"
$ws.Range("R127").Value = "Provincia de Curicó"
$ws.Range("S127").Value = 2500
$ws.Range("T127").Value = 2
